$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bf = New-Object "object[,]" 1,5
$inn = New-Object "object[,]" 1,6

$bf[0,0] = 1.02
$bf[0,1] = 1.036832614200507
$bf[0,2] = 1.04360942655794
$bf[0,3] = 1.049888901486318
$bf[0,4] = 1.056273112988516
$ws.Range("B2:F2").Value = $bf

$inn[0,0] = 1.036448201176315
$inn[0,1] = 1.041938781660735
$inn[0,2] = 1.046382809617151
$inn[0,3] = 1.052644705478561
$inn[0,4] = 1.059011276876586
$inn[0,5] = 1.017916632886698
$ws.Range("I2:N2").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.037765880905275
$bf[0,2] = 1.044323704529178
$bf[0,3] = 1.050777955198726
$bf[0,4] = 1.057173847742474
$ws.Range("B3:F3").Value = $bf

$inn[0,0] = 1.036611809523158
$inn[0,1] = 1.042516258011636
$inn[0,2] = 1.046908410412822
$inn[0,3] = 1.053345897677314
$inn[0,4] = 1.059725396140532
$inn[0,5] = 1.0181089297731
$ws.Range("I3:N3").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.038370150449938
$bf[0,2] = 1.044785823312586
$bf[0,3] = 1.051353905757179
$bf[0,4] = 1.057757233915329
$ws.Range("B4:F4").Value = $bf

$inn[0,0] = 1.036715888985716
$inn[0,1] = 1.042889683197503
$inn[0,2] = 1.047247770740151
$inn[0,3] = 1.053799652200669
$inn[0,4] = 1.060187393988776
$inn[0,5] = 1.018233242534517
$ws.Range("I4:N4").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.03862427609896
$bf[0,2] = 1.044980080508854
$bf[0,3] = 1.05159619510375
$bf[0,4] = 1.0580026197267
$ws.Range("B5:F5").Value = $bf

$inn[0,0] = 1.036759215629298
$inn[0,1] = 1.043046612482485
$inn[0,2] = 1.047390259819355
$inn[0,3] = 1.053990418039532
$inn[0,4] = 1.06038159618306
$inn[0,5] = 1.018285475414527
$ws.Range("I5:N5").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.038666950202338
$bf[0,2] = 1.045012696073362
$bf[0,3] = 1.051636885875081
$bf[0,4] = 1.058043828678095
$ws.Range("B6:F6").Value = $bf

$inn[0,0] = 1.036766465221502
$inn[0,1] = 1.043072958129156
$inn[0,2] = 1.047414173889723
$inn[0,3] = 1.054022448863036
$inn[0,4] = 1.060414202274536
$inn[0,5] = 1.018294243879592
$ws.Range("I6:N6").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.038373545734158
$bf[0,2] = 1.044788419057621
$bf[0,3] = 1.051357142613503
$bf[0,4] = 1.057760512262643
$ws.Range("B7:F7").Value = $bf

$inn[0,0] = 1.036716469602596
$inn[0,1] = 1.042891780326128
$inn[0,2] = 1.047249675387142
$inn[0,3] = 1.053802201195496
$inn[0,4] = 1.060189989014958
$inn[0,5] = 1.018233940584135
$ws.Range("I7:N7").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.037147935970127
$bf[0,2] = 1.043850833054272
$bf[0,3] = 1.050189221044367
$bf[0,4] = 1.056577405594753
$ws.Range("B8:F8").Value = $bf

$inn[0,0] = 1.036503862739687
$inn[0,1] = 1.042133991686692
$inn[0,2] = 1.046560591067584
$inn[0,3] = 1.052881668402462
$inn[0,4] = 1.059252633637205
$inn[0,5] = 1.017981644341732
$ws.Range("I8:N8").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.034991241026021
$bf[0,2] = 1.042198242687991
$bf[0,3] = 1.048136415883399
$bf[0,4] = 1.054496904849935
$ws.Range("B9:F9").Value = $bf

$inn[0,0] = 1.036115576891601
$inn[0,1] = 1.040796874796531
$inn[0,2] = 1.045340732694574
$inn[0,3] = 1.051259904827268
$inn[0,4] = 1.057600300917688
$inn[0,5] = 1.017536193480158
$ws.Range("I9:N9").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.033555515673857
$bf[0,2] = 1.041096307675584
$bf[0,3] = 1.046771476265163
$bf[0,4] = 1.053112875184948
$ws.Range("B10:F10").Value = $bf

$inn[0,0] = 1.035847590228416
$inn[0,1] = 1.039904317392146
$inn[0,2] = 1.044523796018747
$inn[0,3] = 1.050179025972451
$inn[0,4] = 1.056498425394026
$inn[0,5] = 1.017238662391476
$ws.Range("I10:N10").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.032934334970543
$bf[0,2] = 1.040619127078086
$bf[0,3] = 1.046181313412571
$bf[0,4] = 1.052514298401298
$ws.Range("B11:F11").Value = $bf

$inn[0,0] = 1.035729392926006
$inn[0,1] = 1.039517571062288
$inn[0,2] = 1.044169190804944
$inn[0,3] = 1.049711079016798
$inn[0,4] = 1.056021241306863
$inn[0,5] = 1.017109699141921
$ws.Range("I11:N11").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.032703676582799
$bf[0,2] = 1.040441877111642
$bf[0,3] = 1.045962232019209
$bf[0,4] = 1.052292069555515
$ws.Range("B12:F12").Value = $bf

$inn[0,0] = 1.035685165777556
$inn[0,1] = 1.039373877720019
$inn[0,2] = 1.044037345620259
$inn[0,3] = 1.04953727599162
$inn[0,4] = 1.055843985432822
$inn[0,5] = 1.01706177723898
$ws.Range("I12:N12").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.032753150144732
$bf[0,2] = 1.040479897981441
$bf[0,3] = 1.046009219754544
$bf[0,4] = 1.052339733431115
$ws.Range("B13:F13").Value = $bf

$inn[0,0] = 1.035694667274646
$inn[0,1] = 1.03940470215766
$inn[0,2] = 1.044065632680343
$inn[0,3] = 1.04957455670928
$inn[0,4] = 1.055882007781956
$inn[0,5] = 1.0170720575132
$ws.Range("I13:N13").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.032915267124583
$bf[0,2] = 1.040604475610476
$bf[0,3] = 1.046163201388266
$bf[0,4] = 1.052495926658954
$ws.Range("B14:F14").Value = $bf

$inn[0,0] = 1.035725743689478
$inn[0,1] = 1.039505694106188
$inn[0,2] = 1.044158295066325
$inn[0,3] = 1.049696712126184
$inn[0,4] = 1.056006589446246
$inn[0,5] = 1.017105738292217
$ws.Range("I14:N14").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.033015162818014
$bf[0,2] = 1.040681231533897
$bf[0,3] = 1.046258092004354
$bf[0,4] = 1.052592176978796
$ws.Range("B15:F15").Value = $bf

$inn[0,0] = 1.035744848060651
$inn[0,1] = 1.039567913502074
$inn[0,2] = 1.044215370357775
$inn[0,3] = 1.049771977909122
$inn[0,4] = 1.056083347238417
$inn[0,5] = 1.017126487597313
$ws.Range("I15:N15").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.033596751964687
$bf[0,2] = 1.041127975959725
$bf[0,3] = 1.046810661775146
$bf[0,4] = 1.053152615990989
$ws.Range("B16:F16").Value = $bf

$inn[0,0] = 1.03585538920223
$inn[0,1] = 1.03992997901402
$inn[0,2] = 1.044547311844207
$inn[0,3] = 1.050210083872419
$inn[0,4] = 1.056530093274942
$inn[0,5] = 1.01724721854342
$ws.Range("I16:N16").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.033961701456086
$bf[0,2] = 1.0414081984884
$bf[0,3] = 1.047157506777074
$bf[0,4] = 1.053504357631081
$ws.Range("B17:F17").Value = $bf

$inn[0,0] = 1.035924151664955
$inn[0,1] = 1.04015702340867
$inn[0,2] = 1.044755298845896
$inn[0,3] = 1.05048491861538
$inn[0,4] = 1.05681030882926
$inn[0,5] = 1.017322915246056
$ws.Range("I17:N17").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.034174618288883
$bf[0,2] = 1.041571644027425
$bf[0,3] = 1.047359898965082
$bf[0,4] = 1.053709591679058
$ws.Range("B18:F18").Value = $bf

$inn[0,0] = 1.035964051450091
$inn[0,1] = 1.040289429112779
$inn[0,2] = 1.044876530490211
$inn[0,3] = 1.050645232780442
$inn[0,4] = 1.056973747492615
$inn[0,5] = 1.017367055252827
$ws.Range("I18:N18").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.034247225532929
$bf[0,2] = 1.04162737409002
$bf[0,3] = 1.047428923568398
$bf[0,4] = 1.053779582861436
$ws.Range("B19:F19").Value = $bf

$inn[0,0] = 1.035977620910604
$inn[0,1] = 1.040334571689292
$inn[0,2] = 1.04491785312444
$inn[0,3] = 1.050699897039466
$inn[0,4] = 1.057029474741197
$inn[0,5] = 1.017382103706146
$ws.Range("I19:N19").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.033922540860673
$bf[0,2] = 1.041378133616519
$bf[0,3] = 1.047120284957667
$bf[0,4] = 1.053466611923901
$ws.Range("B20:F20").Value = $bf

$inn[0,0] = 1.035916795633046
$inn[0,1] = 1.04013266631908
$inn[0,2] = 1.044732992453651
$inn[0,3] = 1.050455430642483
$inn[0,4] = 1.056780244997636
$inn[0,5] = 1.017314795006761
$ws.Range("I20:N20").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.032867525614746
$bf[0,2] = 1.040567790698038
$bf[0,3] = 1.046117853997349
$bf[0,4] = 1.052449928619461
$ws.Range("B21:F21").Value = $bf

$inn[0,0] = 1.035716601381869
$inn[0,1] = 1.039475955550336
$inn[0,2] = 1.044131011855771
$inn[0,3] = 1.049660740023664
$inn[0,4] = 1.055969903474554
$inn[0,5] = 1.017095820669592
$ws.Range("I21:N21").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.032204634605932
$bf[0,2] = 1.040058274420084
$bf[0,3] = 1.045488346611166
$bf[0,4] = 1.051811332132807
$ws.Range("B22:F22").Value = $bf

$inn[0,0] = 1.035588860519665
$inn[0,1] = 1.039062832959123
$inn[0,2] = 1.043751776847495
$inn[0,3] = 1.049161164078144
$inn[0,4] = 1.055460361288199
$inn[0,5] = 1.016958031852138
$ws.Range("I22:N22").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.0325560037347
$bf[0,2] = 1.040328380186124
$bf[0,3] = 1.045821987807661
$bf[0,4] = 1.052149803753389
$ws.Range("B23:F23").Value = $bf

$inn[0,0] = 1.035656755462967
$inn[0,1] = 1.039281857818598
$inn[0,2] = 1.043952886830726
$inn[0,3] = 1.049425990990794
$inn[0,4] = 1.055730483399617
$inn[0,5] = 1.017031086705884
$ws.Range("I23:N23").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.033940235699931
$bf[0,2] = 1.041391718649018
$bf[0,3] = 1.047137103638205
$bf[0,4] = 1.053483667370858
$ws.Range("B24:F24").Value = $bf

$inn[0,0] = 1.035920120150575
$inn[0,1] = 1.040143672318145
$inn[0,2] = 1.044743072010715
$inn[0,3] = 1.050468754964423
$inn[0,4] = 1.05679382956873
$inn[0,5] = 1.017318464232329
$ws.Range("I24:N24").Value = $inn

$bf[0,0] = 1.02
$bf[0,1] = 1.035548437946586
$bf[0,2] = 1.04262552018734
$bf[0,3] = 1.048666487502889
$bf[0,4] = 1.055034247125942
$ws.Range("B25:F25").Value = $bf

$inn[0,0] = 1.036217570491877
$inn[0,1] = 1.041142758086058
$inn[0,2] = 1.045656752177976
$inn[0,3] = 1.051679122599801
$inn[0,4] = 1.058027530708901
$inn[0,5] = 1.017651454257335
$ws.Range("I25:N25").Value = $inn

